# Apply the changes described by the diff:
# 1. Rename sheet "upper_primary" -> "intermediate level"
# 2. Delete column D (the redundant "Category" column) from the
#    ECE, primary, upper_primary (now intermediate level) and secondary sheets.
#    The Disability sheet keeps its column D untouched.

$wb = $excel.ActiveWorkbook

# 1. Rename the "upper_primary" sheet (updates both the tab name and the
#    text shown in its own header cell, since that cell's text equals the
#    old sheet name).
$wsUpperPrimary = $wb.Worksheets.Item("upper_primary")
$wsUpperPrimary.Range("C1").Value2 = "intermediate level"
$wsUpperPrimary.Name = "intermediate level"

# 2. Remove column D (redundant "Category" column) from the relevant sheets
$sheetsToTrim = @("ECE", "primary", "intermediate level", "secondary")
foreach ($sheetName in $sheetsToTrim) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Columns.Item(4).Delete()
}
